$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#staatz"
$ws.Range("C2").Value = "Staatz"
$ws.Range("D2").ClearContents()

$ws.Range("B3").Value = "#mars"
$ws.Range("C3").Value = "Mars"
$ws.Range("D3").ClearContents()

$ws.Range("B4").Value = "#holland-met-de-3-prov.-voorsigtigheit"
$ws.Range("C4").Value = "Holland met de 3 Prov. Voorsigtigheit"
$ws.Range("D4").ClearContents()

$ws.Range("B5").Value = "#roof"
$ws.Range("C5").Value = "Roof"
$ws.Range("D5").ClearContents()

$ws.Range("B6").Value = "#holland"
$ws.Range("C6").Value = "Holland"
$ws.Range("D6").ClearContents()

$ws.Range("B7").Value = "#de-faam"
$ws.Range("C7").Value = "De Faam"
$ws.Range("D7").ClearContents()

$ws.Range("B8").Value = "#zorgel"
$ws.Range("C8").Value = "Zorgel"
$ws.Range("D8").ClearContents()

$ws.Range("B9").Value = "#vrede"
$ws.Range("C9").Value = "Vrede"
$ws.Range("D9").ClearContents()

$ws.Range("B10").Value = "#voorzig"
$ws.Range("C10").Value = "Voorzig"
$ws.Range("D10").ClearContents()

$ws.Range("B11").Value = "#voorzt:"
$ws.Range("C11").Value = "Voorzt:"
$ws.Range("D11").ClearContents()

$ws.Range("B12").Value = "#holl"
$ws.Range("C12").Value = "Holl"
$ws.Range("D12").ClearContents()

$ws.Range("B13").Value = "#vrank"
$ws.Range("C13").Value = "Vrank"
$ws.Range("D13").ClearContents()

$ws.Range("B14").Value = "#welv"
$ws.Range("C14").Value = "Welv"
$ws.Range("D14").ClearContents()

$ws.Range("B15").Value = "#mart"
$ws.Range("C15").Value = "Mart"
$ws.Range("D15").ClearContents()

$ws.Range("B16").Value = "#weelde"
$ws.Range("C16").Value = "Weelde"
$ws.Range("D16").ClearContents()

$ws.Range("B17").Value = "#spanje"
$ws.Range("C17").Value = "Spanje"
$ws.Range("D17").ClearContents()

$ws.Range("B18").Value = "#faam"
$ws.Range("C18").Value = "Faam"
$ws.Range("D18").ClearContents()

$ws.Range("B19").Value = "#tijd"
$ws.Range("C19").Value = "Tijd"
$ws.Range("D19").ClearContents()

$ws.Range("B20").Value = "#de-provincie-van-holland.-voorzichtigheyd"
$ws.Range("C20").Value = "De Provincie van Holland. Voorzichtigheyd"
$ws.Range("D20").ClearContents()

$ws.Range("B21").Value = "#vrankrijk"
$ws.Range("C21").Value = "Vrankrijk"
$ws.Range("D21").ClearContents()

$ws.Range("B22").Value = "#vreede"
$ws.Range("C22").Value = "Vreede"
$ws.Range("D22").ClearContents()

$ws.Range("B23").Value = "#verraad"
$ws.Range("C23").Value = "Verraad"
$ws.Range("D23").ClearContents()

$ws.Range("B24").Value = "#staatk"
$ws.Range("C24").Value = "Staatk"
$ws.Range("D24").ClearContents()

$ws.Range("B25").Value = "#geweld"
$ws.Range("C25").Value = "Geweld"
$ws.Range("D25").ClearContents()

$ws.Range("B26").Value = "#godv"
$ws.Range("C26").Value = "Godv"
$ws.Range("D26").ClearContents()

$ws.Range("B27").Value = "#voorz"
$ws.Range("C27").Value = "Voorz"
$ws.Range("D27").ClearContents()

$ws.Range("B28").Value = "#godvr"
$ws.Range("C28").Value = "Godvr"
$ws.Range("D28").ClearContents()

$ws.Range("B29").Value = "#voorlig"
$ws.Range("C29").Value = "Voorlig"
$ws.Range("D29").ClearContents()

$ws.Range("B30").Value = "#overv"
$ws.Range("C30").Value = "Overv"
$ws.Range("D30").ClearContents()

$ws.Range("B31").Value = "#zorg"
$ws.Range("C31").Value = "Zorg"
$ws.Range("D31").ClearContents()

$ws.Range("B32").Value = "#overvl"
$ws.Range("C32").Value = "Overvl"
$ws.Range("D32").ClearContents()

$ws.Range("B33").Value = "#twist"
$ws.Range("C33").Value = "Twist"
$ws.Range("D33").ClearContents()

$ws.Range("B34").Value = "#spanje,"
$ws.Range("C34").Value = "Spanje,"
$ws.Range("D34").ClearContents()

$ws.Range("B35").Value = "#holla"
$ws.Range("C35").Value = "Holla"
$ws.Range("D35").ClearContents()
